# "Plantilla Checklist" - fill in example rows (Actividad 1-4) under the
# existing header rows, matching each activity's response type / answers,
# and leave a formatted (underlined) blank marker cell further down the
# sheet (C12) as a template placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Underlined placeholder cell, far below the data (row 12) - create this
# style first so it lands at cellXfs index 6, matching the workbook's
# eventual save order.
$ws.Range("C12").Font.Underline = $true

# The four example-activity rows share one bold style (cellXfs index 7).
$ws.Range("A3:E6").Font.Bold = $true

# Row 3 - Actividad 1 / Opción única
$ws.Range("A3").Value = "Actividad 1"
$ws.Range("B3").Value = "Opción única"
$ws.Range("C3").Value = "si"
$ws.Range("D3").Value = "no"
$ws.Range("E3").Value = "si"

# Row 4 - Actividad 2 / Escala de evaluación
$ws.Range("A4").Value = "Actividad 2"
$ws.Range("B4").Value = "Escala de evaluación"
$ws.Range("C4").Value = "no"
$ws.Range("D4").Value = "no"
$ws.Range("E4").Value = "no"

# Row 5 - Actividad 3 / Opción única
$ws.Range("A5").Value = "Actividad 3"
$ws.Range("B5").Value = "Opción única"
$ws.Range("C5").Value = "si"
$ws.Range("D5").Value = "si"
$ws.Range("E5").Value = "no"

# Row 6 - Actividad 4 / Escala de evaluación
$ws.Range("A6").Value = "Actividad 4"
$ws.Range("B6").Value = "Escala de evaluación"
$ws.Range("C6").Value = "no"
$ws.Range("D6").Value = "si"
$ws.Range("E6").Value = "si"

# Leave the selection where the author left it when saving.
$ws.Range("B6").Select() | Out-Null
